# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

# --- "Forecast Comparison" sheet: MyForecast column adjustments ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D4").Value = 172
$wsForecast.Range("D6").Value = 243

# --- "Summary" sheet: totals/min recalculated ---
# Column B on this sheet stores numbers-as-text (inline string cells), so
# force Text format before assigning to avoid Excel re-inferring a Number type.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "3944"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "2070"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "1074"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "172"
